# Apply the "VENTA MENSUAL" PRESUPUESTO (G column) edit:
# The G2:G38 cells currently hold text values like "$2.300,00 " (shared
# strings). Replace them with real numeric values and apply a currency
# number format, leaving the blank G cells untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

$values = @{
    2  = 2300
    3  = 500
    4  = 500
    6  = 1500
    9  = 500
    11 = 1500
    14 = 1000
    15 = 1500
    16 = 500
    20 = 500
    21 = 500
    25 = 800
    27 = 1500
    31 = 2000
    33 = 500
    34 = 500
    35 = 1500
    36 = 2410
    37 = 1617
    38 = 1000
}

foreach ($row in $values.Keys) {
    $cell = $ws.Cells.Item($row, 7)
    $cell.Value = $values[$row]
    $cell.NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"
}

# Update the sheet's current selection to match the edited range.
$ws.Activate()
$ws.Range("G2:G38").Select()
